$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: B2 date value and C2 status text
$ws.Range("B2").Value = 46000
$ws.Range("C2").Value = "Sent at 2025-12-09 15:37:14"

# Update row 3: B3 date value and C3 status text
$ws.Range("B3").Value = 46000
$ws.Range("C3").Value = "Sent at 2025-12-09 15:37:16"

# Remove row 4 entirely (delete row, shifting cells up)
$ws.Rows("4:4").Delete()
